# TPS Screen Project - add 4/28/2018 info
# Mirrors the existing "4/21/2018" block (rows 1-7) into a new block
# at rows 11-17 with the new task info, and updates the "% Complete"
# values in the first block from hours(10) to fractions (0.1 = 10%).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Copy the formatting (fill/border/font/number format) of the
#    first block (rows 1-7) down onto the new block (rows 11-17),
#    so the new section visually matches the existing one.
# ---------------------------------------------------------------
$ws.Range("A1:G7").Copy() | Out-Null
$ws.Range("A11:G17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Re-apply the percentage style (fill + thin border + 0% format,
# same as F5) onto every "% Complete" cell that must show a percent
# number format: F6, F7 (existing block) and F15/F16/F17 (new block).
$ws.Range("F5").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").PasteSpecial(-4122) | Out-Null
$ws.Range("F16").PasteSpecial(-4122) | Out-Null
$ws.Range("F17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 2) Fix up values in the existing block: "% Complete" is now
#    stored as a fraction (0.1) instead of a raw number (10).
# ---------------------------------------------------------------
$ws.Range("F6").Value = 0.1
$ws.Range("F7").Value = 0.1

# ---------------------------------------------------------------
# 3) Merge the cells for the new block the same way as the first.
# ---------------------------------------------------------------
$ws.Range("E11:G11").Merge()
$ws.Range("A13:G13").Merge()

# ---------------------------------------------------------------
# 4) Fill in the content of the new block (rows 11-17).
#    Values are assigned in the same order the original author
#    entered them in, so new shared-string entries land in the
#    same order as the source workbook.
# ---------------------------------------------------------------
# Header row: new date
$ws.Range("A11").Value = "Date: 04/28/2018"

# Section title row (mirrors "TPS" banner on row 3)
$ws.Range("A13").Value = "TPS"

# Column headers (mirrors row 4)
$ws.Range("A14").Value = "Task"
$ws.Range("B14").Value = "Time Est (Hrs)"
$ws.Range("C14").Value = "Risk"
$ws.Range("D14").Value = "Who"
$ws.Range("E14").Value = "Time Spent (Hrs)"
$ws.Range("F14").Value = "% Complete"
$ws.Range("G14").Value = "Peer Review"

# Row 15: Get touch screen to work
$ws.Range("A15").Value = "Get touch screen to work"
$ws.Range("B15").Value = 0.5
$ws.Range("C15").Value = 4
$ws.Range("E15").Value = 0.5
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = "Yes"

# Row 16: Get touch keyboard to work with browser
$ws.Range("A16").Value = "Get touch keyboard to work with browser"
$ws.Range("B16").Value = 6
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 0.75
$ws.Range("G16").Value = "Yes"

# "Who" column for rows 15/16 (same assignee)
$ws.Range("D15").Value = "Matt"
$ws.Range("D16").Value = "Matt"

# Row 17: Test embedded browser
$ws.Range("A17").Value = "Test embedded browser"
$ws.Range("B17").Value = 10
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = "Jonah/Matt"
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 0.1
$ws.Range("G17").Value = "Yes"

# ---------------------------------------------------------------
# 5) Cosmetic touch-ups matching the source edit: widen a couple of
#    columns to fit the new text, and leave the selection where the
#    author left off (just below the new table).
# ---------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 52.17
$ws.Columns("B").ColumnWidth = 13.3
$ws.Columns("C").ColumnWidth = 5.8
$ws.Columns("D").ColumnWidth = 11.5
$ws.Columns("E").ColumnWidth = 16.1
$ws.Columns("F").ColumnWidth = 11.8
$ws.Columns("G").ColumnWidth = 11.5

$ws.Range("G18").Select() | Out-Null

$wb.Save()
